$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.079427599906921
$ws.Range("B1").Value = 2.793058633804321
$ws.Range("C1").Value = 5.16887092590332
$ws.Range("D1").Value = 2.100213527679443
$ws.Range("E1").Value = 1.173641562461853
